$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# A new September entry ("corporate internet share" @ 2024-09-09 11:15:51) was
# logged; it sorts above the existing row 35, so insert a fresh row there and
# push every following row (through the old last row, 101) down by one.
$ws.Rows.Item(35).Insert()

$ws.Range("R35").Value = "corporate internet share"
$ws.Range("S35").Value = "2024-09-09 11:15:51"
